# ---------------------------------------------------------------------------
# Add 2022-Q1 data:
#   1. The current "总计" sheet (holding the running totals table) is renamed
#      to "2022-Q1" and its body is replaced with the new quarter's detailed
#      fund-holdings table (this keeps the original sheetId).
#   2. A brand new sheet named "总计" is inserted right after it, re-created
#      with the same running-totals table plus a new leading row for 2022-Q1
#      (this gets the next sheetId, exactly like the source diff shows).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$refSheet = $wb.Worksheets.Item("2021-Q4")   # used as a formatting template
$q1Sheet  = $wb.Worksheets.Item("总计")       # becomes "2022-Q1"

# ---- Step 1: repurpose the old "总计" sheet into the new "2022-Q1" sheet ----

# Clear out the previous totals content first.
$q1Sheet.Cells.Clear()
$q1Sheet.Name = "2022-Q1"

$fundHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q1Sheet.Cells.Item(1, $col).Value = $fundHeaders[$col - 2]
}

$fundRows = @(
    @("160916", "大成优选混合(LOF)",           "16.14", "89.35", "4.46", "0.7198", 5),
    @("005123", "南方优享分红灵活配置混合A",     "4.86",  "93.90", "6.58", "0.3198", 6),
    @("481017", "工银量化策略混合",             "7.21",  "89.20", "2.83", "0.2040", 4),
    @("010738", "大成优选升级一年持有期混合A",   "3.79",  "89.02", "5.29", "0.2005", 8),
    @("519013", "海富通风格优势混合",           "3.57",  "89.46", "3.96", "0.1414", 4),
    @("000030", "长城核心优选灵活配置混合",      "2.25",  "91.66", "4.59", "0.1033", 2),
    @("010375", "国金鑫悦经济新动能混合A",       "1.29",  "94.93", "7.65", "0.0987", 3),
    @("001648", "工银瑞信新价值灵活配置混合",     "1.96",  "80.70", "2.19", "0.0429", 3),
    @("005562", "创金合信中证红利低波动指数C",   "1.85",  "94.46", "2.05", "0.0379", 10),
    @("002003", "工银瑞信新机遇灵活配置混合A",   "1.40",  "87.77", "2.32", "0.0325", 4),
    @("011765", "兴银高端制造混合A",             "1.01",  "93.23", "3.17", "0.0320", 2),
    @("002004", "工银瑞信新机遇灵活配置混合C",   "1.29",  "87.77", "2.32", "0.0299", 4),
    @("512890", "华泰柏瑞中证红利低波动ETF",     "1.36",  "99.24", "2.16", "0.0294", 10),
    @("005561", "创金合信中证红利低波动指数A",   "1.22",  "94.46", "2.05", "0.0250", 10),
    @("010376", "国金鑫悦经济新动能混合C",       "0.28",  "94.93", "7.65", "0.0214", 3),
    @("011766", "兴银高端制造混合C",             "0.39",  "93.23", "3.17", "0.0124", 2),
    @("002135", "广发鑫源灵活配置混合A",         "0.49",  "26.98", "2.51", "0.0123", 5),
    @("006587", "南方优享分红灵活配置混合C",     "0.09",  "93.90", "6.58", "0.0059", 6),
    @("010739", "大成优选升级一年持有期混合C",   "0.09",  "89.02", "5.29", "0.0048", 8),
    @("002136", "广发鑫源灵活配置混合C",         "0.00",  "26.98", "2.51", 0,        5)
)

$r = 2
foreach ($row in $fundRows) {
    $q1Sheet.Cells.Item($r, 1).Value = $r - 2

    $q1Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1Sheet.Cells.Item($r, 3).Value = $row[1]
    $q1Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1Sheet.Cells.Item($r, 6).Value = "'" + $row[4]

    if ($row[5] -eq 0) {
        $q1Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q1Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
    }

    $q1Sheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---- formatting: reuse the same look as the other quarterly sheets ----
$refSheet.Range("B1:H1").Copy() | Out-Null
$q1Sheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$refSheet.Range("A2").Copy() | Out-Null
$q1Sheet.Range("A2:A21").PasteSpecial(-4122) | Out-Null

$q1Sheet.Application.CutCopyMode = $false

# ---- Step 2: create the new "总计" sheet after "2022-Q1" -------------------

$totalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$totalSheet.Name = "总计"
$totalSheet.Outline.SummaryRow = 1
$totalSheet.Outline.SummaryColumn = 1

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $totalSheet.Cells.Item(1, $col).Value = $totalHeaders[$col - 2]
}

$totalRows = @(
    @("2022-Q1", 20, 2.07),
    @("2021-Q4", 24, 3.53),
    @("2021-Q3", 8, 0.47),
    @("2021-Q1", 1, 0.04),
    @("2020-Q4", 1, 0)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

$refSheet.Range("B1:D1").Copy() | Out-Null
$totalSheet.Range("B1:D1").PasteSpecial(-4122) | Out-Null

$refSheet.Range("A2").Copy() | Out-Null
$totalSheet.Range("A2:A6").PasteSpecial(-4122) | Out-Null

$totalSheet.Application.CutCopyMode = $false

# Match the page-margin convention used by the other sheets in this workbook.
$totalSheet.PageSetup.LeftMargin = 54
$totalSheet.PageSetup.RightMargin = 54
$totalSheet.PageSetup.TopMargin = 72
$totalSheet.PageSetup.BottomMargin = 72
$totalSheet.PageSetup.HeaderMargin = 36
$totalSheet.PageSetup.FooterMargin = 36

# Keep the originally active sheet selected (creating/renaming sheets shifts
# the active tab to the last-touched sheet otherwise).
$wb.Worksheets.Item(1).Activate()
